$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column C to fit the newly-populated oligo-sequence cells ---
# Target stored width (OOXML "width" units) is 56.41 characters; the engine
# quantizes ColumnWidth to whole pixels, so 55.5 is the closest input that
# rounds to the nearest achievable stored width (56.33).
$ws.Columns.Item(3).ColumnWidth = 55.5

# Fill in newly-added assay columns (C-G) for rows 11-19 (oligo design sequences / dot-bracket structures)
$ws.Range("C11").Value = "CATCACTATCAATCCTACATCTTTTCCTAATCCCAATCAACACC"
$ws.Range("D11").Value = "CATCACTATCAATCCTACATCTTTTCCTATCTTCACACCACTCT"
$ws.Range("E11").Value = "GTAAAGAAGAGTGGTGTGAAGATAGGAAAGGTGTTGATTGGGATTAGGAAAAG"
$ws.Range("F11").Value = "TCTTCACACCACTCTTCTTTACTCCTTCAACTCTCCA"
$ws.Range("G11").Value = "....................((((((((((((((((((((((((+......................((((((((((((((((((((((+.......))))))))))))))))))))))))))))))))))))))))))))))"
$ws.Range("C12").Value = "CATCACTATCAATCCTACATCTTTTCCTAATCCCAATCAACACC"
$ws.Range("D12").Value = "TCTTCACACCACTCTTCTTTACTCCTTCAACTCTCCA"
$ws.Range("E12").Value = "GTAAAGAAGAGTGGTGTGAAGATAGGAAAGGTGTTGATTGGGATTAGGAAAAG"
$ws.Range("F12").Value = "CATCACTATCAATCCTACATCTTTTCCTATCTTCACACCACTCT"
$ws.Range("G12").Value = "....................((((((((((((((((((((((((+((((((((((((((((((((((...............+)))))))))))))))))))))).......))))))))))))))))))))))))"
$ws.Range("C13").Value = "CATCACTATCAATCCTACATCTTTTCCTAATCCCAATCAACACC"
$ws.Range("D13").Value = "TCTTCACACCACTCTTCTTTACTCCTTCAACTCTCCA"
$ws.Range("E13").Value = "GTAAAGAAGAGTGGTGTGAAGATAGGAAAGGTGTTGATTGGGATTAGGAAAAG"
$ws.Range("F13").Value = "TTTCCTAATCCCAATCAACACCTTTCCTA"
$ws.Range("G13").Value = "....................((((((((((((((((((((((((+((((((((((((((((((((((...............+)))))))))))))))))))))).......))))))))))))))))))))))))"
$ws.Range("C14").Value = "CCCATTTCTCTAACTAACCACCCTATACCCTTCTTATCCAACCG"
$ws.Range("D14").Value = "CCCATTTCTCTAACTAACCACCCTATACCAAACTTACATTACCG"
$ws.Range("E14").Value = "AGTGATGCGGTAATGTAAGTTTGGTATAGCGGTTGGATAAGAAGGGTATAGGG"
$ws.Range("F14").Value = "AAACTTACATTACCGCATCACTATCAATCCTACATCT"
$ws.Range("G14").Value = "....................((((((((((((((((((((((((+......................((((((((((((((((((((((+.......))))))))))))))))))))))))))))))))))))))))))))))"
$ws.Range("C15").Value = "CCCATTTCTCTAACTAACCACCCTATACCCTTCTTATCCAACCG"
$ws.Range("D15").Value = "AAACTTACATTACCGCATCACTATCAATCCTACATCT"
$ws.Range("E15").Value = "AGTGATGCGGTAATGTAAGTTTGGTATAGCGGTTGGATAAGAAGGGTATAGGG"
$ws.Range("F15").Value = "CCCATTTCTCTAACTAACCACCCTATACCAAACTTACATTACCG"
$ws.Range("G15").Value = "....................((((((((((((((((((((((((+((((((((((((((((((((((...............+)))))))))))))))))))))).......))))))))))))))))))))))))"
$ws.Range("C16").Value = "CCCATTTCTCTAACTAACCACCCTATACCCTTCTTATCCAACCG"
$ws.Range("D16").Value = "AAACTTACATTACCGCATCACTATCAATCCTACATCT"
$ws.Range("E16").Value = "AGTGATGCGGTAATGTAAGTTTGGTATAGCGGTTGGATAAGAAGGGTATAGGG"
$ws.Range("F16").Value = "CTATACCCTTCTTATCCAACCGCTATACC"
$ws.Range("G16").Value = "....................((((((((((((((((((((((((+((((((((((((((((((((((...............+)))))))))))))))))))))).......))))))))))))))))))))))))"
$ws.Range("C17").Value = "TCTTTACTCCTTCAACTCTCCAAACAACATCCTCACACAAACGC"
$ws.Range("D17").Value = "TCTTTACTCCTTCAACTCTCCAAACAACAATCTTCCCTCCACCG"
$ws.Range("E17").Value = "AAATGGGCGGTGGAGGGAAGATTGTTGTTGCGTTTGTGTGAGGATGTTGTTTG"
$ws.Range("F17").Value = "ATCTTCCCTCCACCGCCCATTTCTCTAACTAACCACC"
$ws.Range("G17").Value = "....................((((((((((((((((((((((((+......................((((((((((((((((((((((+.......))))))))))))))))))))))))))))))))))))))))))))))"
$ws.Range("C18").Value = "TCTTTACTCCTTCAACTCTCCAAACAACATCCTCACACAAACGC"
$ws.Range("D18").Value = "ATCTTCCCTCCACCGCCCATTTCTCTAACTAACCACC"
$ws.Range("E18").Value = "AAATGGGCGGTGGAGGGAAGATTGTTGTTGCGTTTGTGTGAGGATGTTGTTTG"
$ws.Range("F18").Value = "TCTTTACTCCTTCAACTCTCCAAACAACAATCTTCCCTCCACCG"
$ws.Range("G18").Value = "....................((((((((((((((((((((((((+((((((((((((((((((((((...............+)))))))))))))))))))))).......))))))))))))))))))))))))"
$ws.Range("C19").Value = "TCTTTACTCCTTCAACTCTCCAAACAACATCCTCACACAAACGC"
$ws.Range("D19").Value = "ATCTTCCCTCCACCGCCCATTTCTCTAACTAACCACC"
$ws.Range("E19").Value = "AAATGGGCGGTGGAGGGAAGATTGTTGTTGCGTTTGTGTGAGGATGTTGTTTG"
$ws.Range("F19").Value = "AACAACATCCTCACACAAACGCAACAACA"
$ws.Range("G19").Value = "....................((((((((((((((((((((((((+((((((((((((((((((((((...............+)))))))))))))))))))))).......))))))))))))))))))))))))"

# --- Row 16 col E sequence cell uses the smaller "LM Mono Caps 10" 11pt font
#     (matches the font already used on the other oligo-sequence cells) ---
$ws.Range("E16").Font.Name = "LM Mono Caps 10"
$ws.Range("E16").Font.Size = 11

# --- Update the sheet's saved view state: scroll back to column A and move
#     the active selection down to C25 (below the populated table) ---
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C25").Select()
